# Insert a new data row at row 254, pushing existing rows 254-316 down to 255-317.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(254).Insert()

$ws.Range("A254").Value = 3
$ws.Range("B254").Value = "Femacal de La Calera"
$ws.Range("C254").Value = "Coquimbo"
$ws.Range("D254").Value = 44754
$ws.Range("E254").Value = 5
$ws.Range("F254").Value = 100112001
$ws.Range("G254").Value = "Berenjena"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 50
$ws.Range("K254").Value = 8000
$ws.Range("L254").Value = 8000
$ws.Range("M254").Value = 8000
$ws.Range("N254").Value = "$/caja 60 unidades"
$ws.Range("O254").Value = "Región de Arica y Parinacota"
$ws.Range("P254").Value = 133
$ws.Range("Q254").Value = 60
$ws.Range("R254").Value = "Hortaliza"
